# Applies the Thu Jun  8 22:07:33 UTC 2023 "cryptos list" refresh described
# by the OOXML diff: updated Price (D) / Volume(1h) (E) figures for every
# coin row, plus two rank swaps (Stellar<->Filecoin @ rows 30/31 and
# Algorand<->Elrond @ rows 48/49) which also touch the Coin (B) and Link (C)
# columns for those four rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that LOOKS numeric ("1.001", "0.5274", ...) while
# keeping the cell a *text* cell, matching the source XML (t="inlineStr").
# A leading apostrophe forces Excel to store it as text instead of silently
# parsing it into a Number; resetting the style back to "Normal" afterwards
# clears the quotePrefix style bit Excel stamps on the cell so no stray
# style index is left behind.
function Set-TextValue($range, $text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

$ws.Range("D2").Value = '26.608.36'
$ws.Range("E2").Value = '  +1.24%  '
$ws.Range("D3").Value = '1.853.96'
$ws.Range("E3").Value = '  +1.20%  '
Set-TextValue $ws.Range("D4") '1.001'
$ws.Range("E4").Value = '  -0.17%  '
Set-TextValue $ws.Range("D5") '263.66'
$ws.Range("E5").Value = '  +1.66%  '
$ws.Range("E6").Value = '  -0.21%  '
Set-TextValue $ws.Range("D7") '0.5274'
$ws.Range("E7").Value = '  +1.65%  '
Set-TextValue $ws.Range("D8") '0.3241'
$ws.Range("E8").Value = '  +0.88%  '
Set-TextValue $ws.Range("D9") '0.06789'
$ws.Range("E9").Value = '  +1.01%  '
Set-TextValue $ws.Range("D10") '18.96'
$ws.Range("E10").Value = '  +1.11%  '
Set-TextValue $ws.Range("D11") '0.7832'
$ws.Range("E11").Value = '  +2.12%  '
Set-TextValue $ws.Range("D12") '0.07763'
$ws.Range("E12").Value = '  +1.21%  '
$ws.Range("D13").Value = '1.861.59'
$ws.Range("E13").Value = '  -0.47%  '
Set-TextValue $ws.Range("D14") '88.67'
$ws.Range("E14").Value = '  -0.03%  '
Set-TextValue $ws.Range("D15") '5.035'
$ws.Range("E15").Value = '  +0.55%  '
Set-TextValue $ws.Range("D16") '1.001'
$ws.Range("E16").Value = '  -0.23%  '
Set-TextValue $ws.Range("D17") '13.97'
$ws.Range("E17").Value = '  -0.60%  '
$ws.Range("E18").Value = '  -0.21%  '
Set-TextValue $ws.Range("D19") '0.000007956'
$ws.Range("E19").Value = '  +1.51%  '
$ws.Range("D20").Value = '26.634.64'
$ws.Range("E20").Value = '  +1.14%  '
Set-TextValue $ws.Range("D21") '4.635'
$ws.Range("E21").Value = '  +2.40%  '
Set-TextValue $ws.Range("D22") '9.473'
$ws.Range("E22").Value = '  +0.89%  '
Set-TextValue $ws.Range("D23") '6.011'
$ws.Range("E23").Value = '  +2.42%  '
Set-TextValue $ws.Range("D24") '143.71'
$ws.Range("E24").Value = '  -0.63%  '
$ws.Range("E25").Value = '  -5.82%  '
Set-TextValue $ws.Range("D26") '1.678'
$ws.Range("E26").Value = '  +0.83%  '
Set-TextValue $ws.Range("D27") '17.02'
$ws.Range("E27").Value = '  +0.72%  '
Set-TextValue $ws.Range("D28") '111.62'
$ws.Range("E28").Value = '  +0.78%  '
Set-TextValue $ws.Range("D29") '4.183'
$ws.Range("E29").Value = '  +0.56%  '
$ws.Range("B30").Value = 'Stellar'
$ws.Range("C30").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range("D30") '0.08717'
$ws.Range("E30").Value = '  +0.00%  '
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range("D31") '4.102'
$ws.Range("E31").Value = '  +0.36%  '
Set-TextValue $ws.Range("D32") '0.04863'
$ws.Range("E32").Value = '  +0.95%  '
Set-TextValue $ws.Range("D33") '1.132'
$ws.Range("E33").Value = '  +0.64%  '
Set-TextValue $ws.Range("D34") '0.7209'
$ws.Range("E34").Value = '  +6.09%  '
Set-TextValue $ws.Range("D35") '2.867'
$ws.Range("E35").Value = '  +0.36%  '
Set-TextValue $ws.Range("D36") '3.114'
$ws.Range("E36").Value = '  +0.91%  '
Set-TextValue $ws.Range("D37") '2.263'
$ws.Range("E37").Value = '  +2.69%  '
Set-TextValue $ws.Range("D38") '0.01790'
$ws.Range("E38").Value = '  +0.76%  '
Set-TextValue $ws.Range("D39") '0.4877'
$ws.Range("E39").Value = '  -0.42%  '
Set-TextValue $ws.Range("D40") '0.9009'
$ws.Range("E40").Value = '  +0.34%  '
Set-TextValue $ws.Range("D41") '111.25'
$ws.Range("E41").Value = '  -0.17%  '
Set-TextValue $ws.Range("D42") '5.969'
$ws.Range("E42").Value = '  -2.77%  '
$ws.Range("E43").Value = '  -0.24%  '
Set-TextValue $ws.Range("D44") '7.678'
$ws.Range("E44").Value = '  -0.11%  '
Set-TextValue $ws.Range("D45") '0.4196'
$ws.Range("E45").Value = '  +0.56%  '
Set-TextValue $ws.Range("D46") '0.05886'
$ws.Range("E46").Value = '  +0.09%  '
Set-TextValue $ws.Range("D47") '8.997'
$ws.Range("E47").Value = '  -0.05%  '
$ws.Range("B48").Value = 'Algorand'
$ws.Range("C48").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range("D48") '0.1236'
$ws.Range("E48").Value = '  -1.26%  '
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextValue $ws.Range("D49") '35.08'
$ws.Range("E49").Value = '  -0.35%  '
Set-TextValue $ws.Range("D50") '0.8897'
$ws.Range("E50").Value = '  +3.72%  '
Set-TextValue $ws.Range("D51") '59.95'
$ws.Range("E51").Value = '  +1.52%  '
